# Update "want to go" counts (column F) on both the "展览" and "全部类型"
# sheets, which hold duplicated data. Each sheet gets the same set of updates.

$wb = $excel.ActiveWorkbook

# Map of row number -> new value for column F.
$updates = @{
    3  = 8054
    4  = 134
    8  = 137
    9  = 137
    10 = 185
    13 = 164
    14 = 2144
    17 = 16
    20 = 71
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value2 = $updates[$row]
    }
}
